# Generate Report for Handoff
# Replaces the old localization run's file id (f5c79018-...) with the new
# one (fb79a795-...) across all three sheets, refreshes the xliff hash in
# the handback file names, and bumps the handoff/generate timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "f5c79018-1aa1-45e5-a063-f97279efc00a"
$newId = "fb79a795-905e-4de4-91c2-ec580ac9116b"

$oldHash = "7cf58021c14341b91df661944a698e70b73868be"
$newHash = "1e260ded4b1b56e70c10e3a6cc08507ffaaa9186"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newId.md"

$hl = $ws.Range("B2").Hyperlinks.Item(1)
$hl.TextToDisplay = "e2e\$newId.md"

$ws.Range("G2").Value = "2016-08-19 00:55:15"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$hl = $ws.Range("A2").Hyperlinks.Item(1)
$hl.TextToDisplay = "$newId.md"

$ws.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-19 00:55:10"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$hl = $ws.Range("A2").Hyperlinks.Item(1)
$hl.TextToDisplay = "$newId.md"

$ws.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-19 00:55:15"
